# Insert a new row at row 770 (shifting existing rows 770-811 down to 771-812)
# and populate it with the new data point: 2026/02/04, 水, 16, 201

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 770 (Excel default shifts cells down).
$ws.Rows(770).Insert()

# Column A holds a date-looking string ("2026/02/04") that must stay literal text,
# exactly like the rest of column A in this sheet (inline/shared text, not a real
# date serial). Typing it straight into Value would make Excel auto-convert it to a
# date. Using a text formula first, then converting it to a plain value via
# copy/paste-special keeps it as literal text without leaving a formula behind and
# without introducing a new (unused) cell style.
$ws.Range("A770").Formula = "=""2026/02/04"""
$ws.Range("A770").Copy()
$ws.Range("A770").PasteSpecial(-4163)   # xlPasteValues

$ws.Range("B770").Value = "水"
$ws.Range("C770").Value = 16
$ws.Range("D770").Value = 201
